# Updates the "cryptos" worksheet with refreshed price / 1h-volume / hour
# values, mirroring the GitHub Actions scrape commit made on
# Thu Jan 26 22:10:17 UTC 2023.
#
# The sheet stores every data value (Price, Volume(1h), Hora, ...) as plain
# text (inline/shared string), never as a numeric or percentage cell. Excel's
# COM layer will happily auto-coerce a string like "304.73" or "-1.47%" into
# a genuine number/percentage when assigned directly, which would alter the
# cell's type and add number-format styling that should not be there. To
# keep the values as literal text (matching the original authoring), each
# cell is temporarily switched to the "Text" number format before the value
# is written, and then restored to the workbook's default "Normal" style so
# no stray style index is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellUpdates = @(
    @("D2", "304.73"),
    @("E2", "-1.47%"),
    @("G2", "22"),
    @("D3", "35.90"),
    @("E3", "-0.93%"),
    @("G3", "22"),
    @("D4", "4.985"),
    @("E4", "-1.85%"),
    @("G4", "22"),
    @("D5", "0.08100"),
    @("E5", "-1.41%"),
    @("G5", "22"),
    @("D6", "1.904"),
    @("E6", "-5.40%"),
    @("G6", "22"),
    @("D7", "4.147"),
    @("E7", "1.53%"),
    @("G7", "22"),
    @("D8", "7.879"),
    @("E8", "0.13%"),
    @("G8", "22"),
    @("D9", "0.9300"),
    @("E9", "-0.32%"),
    @("G9", "22"),
    @("D10", "0.1287"),
    @("E10", "-14.45%"),
    @("G10", "22"),
    @("D11", "0.1901"),
    @("E11", "-1.51%"),
    @("G11", "22"),
    @("D12", "0.09208"),
    @("E12", "1.63%"),
    @("G12", "22"),
    @("D13", "0.03516"),
    @("E13", "0.68%"),
    @("G13", "22"),
    @("D14", "0.09911"),
    @("E14", "0.28%"),
    @("G14", "22"),
    @("D15", "0.001431"),
    @("E15", "-0.45%"),
    @("G15", "22"),
    @("D16", "0.006616"),
    @("E16", "16.42%"),
    @("G16", "22"),
    @("D17", "3.605"),
    @("E17", "1.34%"),
    @("G17", "22"),
    @("D18", "3.106"),
    @("E18", "1.86%"),
    @("G18", "22"),
    @("D19", "0.3451"),
    @("E19", "0.62%"),
    @("G19", "22"),
    @("D20", "5.241"),
    @("E20", "4.15%"),
    @("G20", "22"),
    @("E21", "0.42%"),
    @("G21", "22"),
    @("D22", "0.2531"),
    @("E22", "6.75%"),
    @("G22", "22"),
    @("D23", "0.04418"),
    @("E23", "-1.41%"),
    @("G23", "22"),
    @("E24", "3.34%"),
    @("G24", "22"),
    @("D25", "0.004710"),
    @("E25", "-4.19%"),
    @("G25", "22"),
    @("D26", "0.0001301"),
    @("E26", "6.75%"),
    @("G26", "22"),
    @("D27", "0.0003130"),
    @("E27", "-28.70%"),
    @("G27", "22"),
    @("G28", "22"),
    @("G29", "22"),
    @("G30", "22"),
    @("G31", "22"),
    @("G32", "22"),
    @("G33", "22"),
    @("G34", "22"),
    @("G35", "22"),
    @("G36", "22"),
    @("G37", "22"),
    @("G38", "22"),
    @("D39", "0.01951"),
    @("E39", "-3.44%"),
    @("G39", "22"),
    @("D40", "0.05213"),
    @("E40", "7.33%"),
    @("G40", "22"),
    @("D41", "0.007558"),
    @("E41", "0.70%"),
    @("G41", "22"),
    @("D42", "0.01017"),
    @("E42", "-0.74%"),
    @("G42", "22"),
    @("D43", "0.1371"),
    @("E43", "-0.43%"),
    @("G43", "22"),
    @("D44", "0.002102"),
    @("E44", "2.29%"),
    @("G44", "22"),
    @("D45", "0.01064"),
    @("E45", "-3.19%"),
    @("G45", "22"),
    @("D46", "0.00006362"),
    @("E46", "4.23%"),
    @("G46", "22"),
    @("E47", "1.24%"),
    @("G47", "22"),
    @("G48", "22"),
    @("D49", "0.001660"),
    @("E49", "41.06%"),
    @("G49", "22"),
    @("D50", "0.00002101"),
    @("E50", "1.24%"),
    @("G50", "22"),
    @("D51", "0.0002001"),
    @("E51", "1.24%"),
    @("G51", "22")
)

foreach ($entry in $cellUpdates) {
    $addr = $entry[0]
    $val = $entry[1]
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}
